# ------------------------------------------------------------------
# Applies the "Second Edit in word file" commit to HTML-Assignment.docx
#   1. Re-orders the opening paragraph: moves the "Write the HTML for
#      the following screen." run to sit right after "Milind Editing---",
#      and adds a brand new list paragraph "Milind Editing 2" (carrying
#      the _GoBack bookmark) right after it.
#   2. Stamps a <w:lastRenderedPageBreak/> in front of the "Apples" run
#      that belongs to the "Roman numbers list:" bullet.
#   3. Adds <w:bCs w:val="0"/> to the "Username:" and "Password:" runs.
#   4. Renumbers the PBrush OLE object's ObjectID.
#   5. Stamps a <w:lastRenderedPageBreak/> in front of the closing
#      "Develop your personal website ..." paragraph.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-ParaText($para) {
    return $para.Range.Text
}

# Find the (1-based) paragraph index whose own text matches $selfLike,
# optionally also requiring the *previous* paragraph's text to match
# $prevLike. Returns -1 if nothing matches.
function Find-Paragraph($doc, $selfLike, $prevLike) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cur = Get-ParaText $doc.Paragraphs($i)
        if ($cur -like $selfLike) {
            if ($prevLike -ne $null) {
                if ($i -le 1) { continue }
                $prev = Get-ParaText $doc.Paragraphs($i - 1)
                if (-not ($prev -like $prevLike)) { continue }
            }
            return $i
        }
    }
    return -1
}

# --------------------------------------------------------------
# 1) Opening paragraph: reorder runs + insert new "Milind Editing 2"
#    paragraph carrying the _GoBack bookmark.
# --------------------------------------------------------------
$openIdx = Find-Paragraph $d "Milind Editing---*" $null
$openRange = $d.Paragraphs($openIdx).Range

$openXml = @'
<w:p w:rsidR="00AF4B4C" w:rsidRPr="00422554" w:rsidRDefault="00F642D2" w:rsidP="00422554"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve">Milind Editing---  </w:t></w:r><w:r w:rsidR="00AF4B4C" w:rsidRPr="00422554"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:t>Write the HTML for the following screen.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:t>Milind Editing 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$openRange.InsertXML($openXml)

# --------------------------------------------------------------
# 2) "Apples" under "Roman numbers list:" gains a lastRenderedPageBreak.
# --------------------------------------------------------------
$applesIdx = Find-Paragraph $d "Apples*" "Roman numbers list:*"
$applesRange = $d.Paragraphs($applesIdx).Range

$applesXml = @'
<w:p w:rsidR="00AF4B4C" w:rsidRPr="00AF4B4C" w:rsidRDefault="00AF4B4C" w:rsidP="00AF4B4C"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r w:rsidRPr="00AF4B4C"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:lastRenderedPageBreak/><w:t>Apples</w:t></w:r></w:p>
'@
$applesRange.InsertXML($applesXml)

# --------------------------------------------------------------
# 3) "Username: " / "Password: " runs gain <w:bCs w:val="0"/>.
# --------------------------------------------------------------
$credIdx = Find-Paragraph $d "Username:*Password:*" $null
$credRange = $d.Paragraphs($credIdx).Range

$credXml = @'
<w:p w:rsidR="00AF4B4C" w:rsidRPr="00AF4B4C" w:rsidRDefault="00AF4B4C" w:rsidP="00AF4B4C"><w:pPr><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r w:rsidRPr="00AF4B4C"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:bCs w:val="0"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve">Username: </w:t></w:r><w:r w:rsidRPr="00AF4B4C"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:object w:dxaOrig="225" w:dyaOrig="225"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1032" type="#_x0000_t75" style="width:49.5pt;height:18pt" o:ole=""><v:imagedata r:id="rId10" o:title=""/></v:shape><w:control r:id="rId11" w:name="DefaultOcxName" w:shapeid="_x0000_i1032"/></w:object></w:r><w:r w:rsidRPr="00AF4B4C"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:bCs w:val="0"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:br/><w:t xml:space="preserve">Password: </w:t></w:r><w:r w:rsidRPr="00AF4B4C"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:object w:dxaOrig="225" w:dyaOrig="225"><v:shape id="_x0000_i1036" type="#_x0000_t75" style="width:49.5pt;height:18pt" o:ole=""><v:imagedata r:id="rId12" o:title=""/></v:shape><w:control r:id="rId13" w:name="DefaultOcxName1" w:shapeid="_x0000_i1036"/></w:object></w:r></w:p>
'@
$credRange.InsertXML($credXml)

# --------------------------------------------------------------
# 4) PBrush OLE object ObjectID renumbered.
# --------------------------------------------------------------
$oleTableIdx = -1
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $txt = $d.Tables($i).Range.Text
    if ($txt -eq "") {
        $oleTableIdx = $i
    }
}
$oleParaRange = $d.Tables($oleTableIdx).Range.Paragraphs(1).Range

$oleXml = @'
<w:p w:rsidR="00AF4B4C" w:rsidRPr="00AF4B4C" w:rsidRDefault="00AF4B4C" w:rsidP="00AF4B4C"><w:pPr><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r w:rsidRPr="00AF4B4C"><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/><w:lang w:val="en-IN"/></w:rPr><w:object w:dxaOrig="3195" w:dyaOrig="1845"><v:shape id="_x0000_i1029" type="#_x0000_t75" style="width:159.75pt;height:92.25pt" o:ole=""><v:imagedata r:id="rId14" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="PBrush" ShapeID="_x0000_i1029" DrawAspect="Content" ObjectID="_1533020563" r:id="rId15"/></w:object></w:r></w:p>
'@
$oleParaRange.InsertXML($oleXml)

# --------------------------------------------------------------
# 5) Closing "Develop your personal website ..." paragraph gains a
#    lastRenderedPageBreak.
# --------------------------------------------------------------
$personalIdx = Find-Paragraph $d "Develop your personal website*" $null
$personalRange = $d.Paragraphs($personalIdx).Range

$personalXml = @'
<w:p w:rsidR="00977947" w:rsidRPr="00977947" w:rsidRDefault="00977947" w:rsidP="00977947"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Candara" w:hAnsi="Candara"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Develop your personal website that provides information about your hobbies, contact details, address, photos, technical skills etc.</w:t></w:r></w:p>
'@
$personalRange.InsertXML($personalXml)

Write-Output "All edits applied."
